# Regenerate merged AHB files
# - Rename the diff-source/diff-target header columns from the generic
#   "_old" / "_new" suffixes to the concrete version tags "_FV2410" / "_FV2504"
# - Freeze the header row
# - Turn the used range into a native Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
    $ws.Range($newSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2504"
}

# Column K1 ("diff") stays untouched.

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into a native Excel table.
$dataRange = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
